$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (Ameliano vs Olimpia Asuncion) - updated odds
$ws.Cells.Item(5, 7).Value2 = 3.6  # G5
$ws.Cells.Item(5, 9).Value2 = 2.25  # I5
$ws.Cells.Item(5, 10).Value2 = 4.5  # J5
$ws.Cells.Item(5, 11).Value2 = 1.8  # K5
$ws.Cells.Item(5, 15).Value2 = 1.67  # O5
$ws.Cells.Item(5, 16).Value2 = 2.1  # P5
$ws.Cells.Item(5, 17).Value2 = 3.4  # Q5
$ws.Cells.Item(5, 18).Value2 = 1.33  # R5
$ws.Cells.Item(5, 19).Value2 = 1.73  # S5
$ws.Cells.Item(5, 20).Value2 = 2  # T5
$ws.Cells.Item(5, 21).Value2 = 2.5  # U5
$ws.Cells.Item(5, 22).Value2 = 1.5  # V5
$ws.Cells.Item(5, 23).Value2 = 7  # W5
$ws.Cells.Item(5, 24).Value2 = 17  # X5
$ws.Cells.Item(5, 31).Value2 = 23  # AE5
$ws.Cells.Item(5, 32).Value2 = 101  # AF5
$ws.Cells.Item(5, 34).Value2 = 5  # AH5
$ws.Cells.Item(5, 35).Value2 = 9  # AI5
$ws.Cells.Item(5, 37).Value2 = 21  # AK5
$ws.Cells.Item(5, 39).Value2 = 51  # AM5
$ws.Cells.Item(5, 40).Value2 = 5.5  # AN5
$ws.Cells.Item(5, 41).Value2 = 23  # AO5
$ws.Cells.Item(5, 46).Value2 = 2  # AT5
$ws.Cells.Item(5, 48).Value2 = 101  # AV5

# Row 6: now River Plate-Penarol (Uruguay) is replaced by Houston Dynamo vs
# Seattle Sounders (USA - MLS), i.e. the match previously on row 7 shifts up
# (with some odds refreshed) and row 7 becomes a brand-new fixture.
$ws.Cells.Item(6, 1).Value2 = "SKm9QS13"  # A6
$ws.Cells.Item(6, 3).Value2 = "20:30"  # C6
$ws.Cells.Item(6, 4).Value2 = "USA - MLS"  # D6
$ws.Cells.Item(6, 5).Value2 = "Houston Dynamo"  # E6
$ws.Cells.Item(6, 6).Value2 = "Seattle Sounders"  # F6
$ws.Cells.Item(6, 7).Value2 = 2.2  # G6
$ws.Cells.Item(6, 8).Value2 = 3.2  # H6
$ws.Cells.Item(6, 9).Value2 = 3.5  # I6
$ws.Cells.Item(6, 10).Value2 = 2.88  # J6
$ws.Cells.Item(6, 11).Value2 = 2.05  # K6
$ws.Cells.Item(6, 12).Value2 = 4  # L6
$ws.Cells.Item(6, 13).Value2 = 1.08  # M6
$ws.Cells.Item(6, 14).Value2 = 8  # N6
$ws.Cells.Item(6, 15).Value2 = 1.4  # O6
$ws.Cells.Item(6, 16).Value2 = 2.75  # P6
$ws.Cells.Item(6, 17).Value2 = 2.25  # Q6
$ws.Cells.Item(6, 18).Value2 = 1.62  # R6
$ws.Cells.Item(6, 19).Value2 = 1.5  # S6
$ws.Cells.Item(6, 20).Value2 = 2.5  # T6
$ws.Cells.Item(6, 21).Value2 = 1.95  # U6
$ws.Cells.Item(6, 22).Value2 = 1.8  # V6
$ws.Cells.Item(6, 23).Value2 = 6.5  # W6
$ws.Cells.Item(6, 24).Value2 = 9.5  # X6
$ws.Cells.Item(6, 25).Value2 = 9.5  # Y6
$ws.Cells.Item(6, 26).Value2 = 21  # Z6
$ws.Cells.Item(6, 27).Value2 = 19  # AA6
$ws.Cells.Item(6, 28).Value2 = 34  # AB6
$ws.Cells.Item(6, 29).Value2 = 8  # AC6
$ws.Cells.Item(6, 30).Value2 = 6  # AD6
$ws.Cells.Item(6, 31).Value2 = 17  # AE6
$ws.Cells.Item(6, 32).Value2 = 51  # AF6
$ws.Cells.Item(6, 33).Value2 = 351  # AG6
$ws.Cells.Item(6, 34).Value2 = 9  # AH6
$ws.Cells.Item(6, 35).Value2 = 17  # AI6
$ws.Cells.Item(6, 36).Value2 = 13  # AJ6
$ws.Cells.Item(6, 37).Value2 = 41  # AK6
$ws.Cells.Item(6, 38).Value2 = 29  # AL6
$ws.Cells.Item(6, 39).Value2 = 41  # AM6
$ws.Cells.Item(6, 40).Value2 = 4  # AN6
$ws.Cells.Item(6, 41).Value2 = 12  # AO6
$ws.Cells.Item(6, 42).Value2 = 23  # AP6
$ws.Cells.Item(6, 43).Value2 = 41  # AQ6
$ws.Cells.Item(6, 44).Value2 = 67  # AR6
$ws.Cells.Item(6, 45).Value2 = 201  # AS6
$ws.Cells.Item(6, 46).Value2 = 2.5  # AT6
$ws.Cells.Item(6, 47).Value2 = 8.5  # AU6
$ws.Cells.Item(6, 49).Value2 = 5.5  # AW6
$ws.Cells.Item(6, 50).Value2 = 21  # AX6
$ws.Cells.Item(6, 51).Value2 = 29  # AY6
$ws.Cells.Item(6, 52).Value2 = 67  # AZ6
$ws.Cells.Item(6, 53).Value2 = 101  # BA6
$ws.Cells.Item(6, 54).Value2 = 251  # BB6
$ws.Cells.Item(6, 55).Value2 = 126  # BC6
$ws.Cells.Item(6, 56).Value2 = 151  # BD6

# Row 7: new fixture, New Mexico vs Phoenix Rising (USA - USL Championship)
$ws.Cells.Item(7, 1).Value2 = "fygoJ8Vi"  # A7
$ws.Cells.Item(7, 3).Value2 = "22:00"  # C7
$ws.Cells.Item(7, 4).Value2 = "USA - USL CHAMPIONSHIP"  # D7
$ws.Cells.Item(7, 5).Value2 = "New Mexico"  # E7
$ws.Cells.Item(7, 6).Value2 = "Phoenix Rising"  # F7
$ws.Cells.Item(7, 7).Value2 = 2.07  # G7
$ws.Cells.Item(7, 8).Value2 = 3.2  # H7
$ws.Cells.Item(7, 9).Value2 = 3.25  # I7
$ws.Cells.Item(7, 10).Value2 = 2.67  # J7
$ws.Cells.Item(7, 11).Value2 = 2.1  # K7
$ws.Cells.Item(7, 12).Value2 = 3.85  # L7
$ws.Cells.Item(7, 14).Value2 = 6.6  # N7
$ws.Cells.Item(7, 15).Value2 = 1.36  # O7
$ws.Cells.Item(7, 16).Value2 = 2.9  # P7
$ws.Cells.Item(7, 17).Value2 = 2.07  # Q7
$ws.Cells.Item(7, 18).Value2 = 1.7  # R7
$ws.Cells.Item(7, 19).Value2 = 1.4  # S7
$ws.Cells.Item(7, 20).Value2 = 2.7  # T7
$ws.Cells.Item(7, 21).Value2 = 1.85  # U7
$ws.Cells.Item(7, 22).Value2 = 1.85  # V7
$ws.Cells.Item(7, 23).Value2 = 6.9  # W7
$ws.Cells.Item(7, 24).Value2 = 9.75  # X7
$ws.Cells.Item(7, 25).Value2 = 8.75  # Y7
$ws.Cells.Item(7, 26).Value2 = 19.5  # Z7
$ws.Cells.Item(7, 27).Value2 = 18  # AA7
$ws.Cells.Item(7, 28).Value2 = 30  # AB7
$ws.Cells.Item(7, 29).Value2 = 6.6  # AC7
$ws.Cells.Item(7, 30).Value2 = 6.2  # AD7
$ws.Cells.Item(7, 32).Value2 = 75  # AF7
$ws.Cells.Item(7, 33).Value2 = 700  # AG7
$ws.Cells.Item(7, 36).Value2 = 11.75  # AJ7
$ws.Cells.Item(7, 37).Value2 = 45  # AK7
$ws.Cells.Item(7, 38).Value2 = 32  # AL7
$ws.Cells.Item(7, 39).Value2 = 40  # AM7
$ws.Cells.Item(7, 41).Value2 = 10.75  # AO7
$ws.Cells.Item(7, 42).Value2 = 19.5  # AP7
$ws.Cells.Item(7, 43).Value2 = 40  # AQ7
$ws.Cells.Item(7, 44).Value2 = 75  # AR7
$ws.Cells.Item(7, 45).Value2 = 250  # AS7
$ws.Cells.Item(7, 46).Value2 = 2.7  # AT7
$ws.Cells.Item(7, 47).Value2 = 7.1  # AU7
$ws.Cells.Item(7, 48).Value2 = 65  # AV7
$ws.Cells.Item(7, 49).Value2 = 5.2  # AW7
$ws.Cells.Item(7, 50).Value2 = 18.5  # AX7
$ws.Cells.Item(7, 51).Value2 = 25  # AY7
$ws.Cells.Item(7, 52).Value2 = 90  # AZ7
$ws.Cells.Item(7, 53).Value2 = 120  # BA7
$ws.Cells.Item(7, 54).Value2 = 350  # BB7
$ws.Cells.Item(7, 55).Value2 = 51  # BC7
$ws.Cells.Item(7, 56).Value2 = 51  # BD7
